# Updates cryptos list values (prices / hourly volume %) per the
# "Updated cryptos list on Sat Jun 15 20:16:41 UTC 2024 with GitHub Actions"
# commit. Also swaps the EnergySwap / dogwifhat rows (49 <-> 51) back to
# their refreshed ranking order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value that *looks* numeric (e.g. "606.77", "0.998") as
# plain text, matching the source data's inline-string cells, instead of
# letting Excel auto-coerce it into a floating point number. We flip the
# cell to Text format, assign the value, then clear the format again so the
# cell is left without any special number formatting (just like the
# original cells, which carry no style).
function Set-TextValue($addr, $val) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.ClearFormats()
}

$ws.Range("D2").Value = "66.122.60"
$ws.Range("E2").Value = "  +0.77%  "
$ws.Range("D3").Value = "3.552.37"
$ws.Range("E3").Value = "  +4.10%  "
$ws.Range("E4").Value = "  +0.09%  "
Set-TextValue "D5" "606.77"
$ws.Range("E5").Value = "  +1.63%  "
Set-TextValue "D6" "144.73"
$ws.Range("E6").Value = "  +1.85%  "
$ws.Range("D7").Value = "3.551.93"
$ws.Range("E7").Value = "  +4.08%  "
$ws.Range("E8").Value = "  +0.14%  "
Set-TextValue "D9" "0.495"
$ws.Range("E9").Value = "  +5.23%  "
$ws.Range("E10").Value = "  +2.12%  "
Set-TextValue "D11" "7.97"
$ws.Range("E11").Value = "  +0.36%  "
$ws.Range("E12").Value = "  +1.95%  "
$ws.Range("D13").Value = "4.157.05"
$ws.Range("E13").Value = "  +4.29%  "
$ws.Range("E14").Value = "  +3.28%  "
Set-TextValue "D15" "30.08"
$ws.Range("E15").Value = "  +1.80%  "
$ws.Range("D16").Value = "3.555.65"
$ws.Range("E16").Value = "  +4.40%  "
$ws.Range("D17").Value = "66.234.49"
$ws.Range("E17").Value = "  +0.95%  "
$ws.Range("E19").Value = "  +9.79%  "
$ws.Range("E20").Value = "  +1.32%  "
Set-TextValue "D21" "14.88"
$ws.Range("E21").Value = "  +2.47%  "
Set-TextValue "D22" "430.34"
$ws.Range("E22").Value = "  +3.85%  "
Set-TextValue "D23" "0.611"
$ws.Range("E23").Value = "  +5.92%  "
Set-TextValue "D24" "79.16"
$ws.Range("E24").Value = "  +2.64%  "
$ws.Range("D25").Value = "3.696.52"
$ws.Range("E25").Value = "  +4.27%  "
$ws.Range("E27").Value = "  +8.52%  "
Set-TextValue "D28" "2.52"
$ws.Range("E28").Value = "  +4.40%  "
Set-TextValue "D29" "7.97"
$ws.Range("E29").Value = "  +2.25%  "
Set-TextValue "D30" "9.09"
$ws.Range("E30").Value = "  -1.68%  "
Set-TextValue "D31" "0.998"
$ws.Range("E31").Value = "  -0.15%  "
$ws.Range("E32").Value = "  +0.86%  "
$ws.Range("E33").Value = "  +4.17%  "
$ws.Range("D34").Value = "3.549.67"
$ws.Range("E34").Value = "  +4.22%  "
Set-TextValue "D35" "0.155"
$ws.Range("E35").Value = "  -3.53%  "
$ws.Range("E37").Value = "  +4.03%  "
$ws.Range("E38").Value = "  +5.33%  "
Set-TextValue "D39" "5.62"
$ws.Range("E39").Value = "  +2.13%  "
$ws.Range("E40").Value = "  +0.11%  "
Set-TextValue "D41" "174.07"
$ws.Range("E41").Value = "  +3.08%  "
Set-TextValue "D42" "0.0852"
$ws.Range("E42").Value = "  -0.09%  "
$ws.Range("E43").Value = "  +3.74%  "
Set-TextValue "D44" "0.892"
$ws.Range("E44").Value = "  +2.42%  "
$ws.Range("E45").Value = "  +1.46%  "
$ws.Range("E46").Value = "  +1.68%  "
Set-TextValue "D47" "1.21"
$ws.Range("E47").Value = "  +2.32%  "
$ws.Range("E48").Value = "  -1.70%  "
$ws.Range("B49").Value = "dogwifhat"
$ws.Range("C49").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-TextValue "D49" "2.35"
$ws.Range("E49").Value = "  +4.27%  "
Set-TextValue "D50" "7.13"
$ws.Range("E50").Value = "  +1.40%  "
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue "D51" "23.46"
$ws.Range("E51").Value = "  +16.07%  "

Write-Output "cryptos list updated"
